$wb = $excel.ActiveWorkbook

# --- Sheet "Parts": add two new computer parts ------------------------
$wsParts = $wb.Worksheets.Item("Parts")
$wsParts.Range("A3").Value = "Intel Core i5-13400F"
$wsParts.Range("B3").Value = "Processor"
$wsParts.Range("A4").Value = "B650 AORUS"
$wsParts.Range("B4").Value = "Motherboard"

# --- Sheet "Webshop": remove the Emag and Alza rows, keep only Pcx ----
$wsShop = $wb.Worksheets.Item("Webshop")
$wsShop.Hyperlinks.Delete()
$wsShop.Rows.Item(2).Delete()
$wsShop.Rows.Item(2).Delete()
$wsShop.Hyperlinks.Add($wsShop.Range("B2"), "http://pcx.hu/")
$wsShop.Range("B2").Style = "Hyperlink"

# --- restore the active selections shown in the saved workbook --------
# (set the non-active sheet's selection first - selecting a range on a
# sheet implicitly activates it, so the sheet meant to stay on top, i.e.
# Parts, must be activated/selected last)
$wsShop.Range("B8").Select() | Out-Null

$wsParts.Activate() | Out-Null
$wsParts.Range("A9").Select() | Out-Null
